# Edit script: applies the author's changes to slide with sldId=265
# (position 11 in the slide list):
#   1. Recolor the "-0.037" delta text from green (7BC043) to red (FF0000).
#   2. Add a new highlighted conclusion textbox ("Shape 13") summarizing
#      that the Tuning model has the best metrics.

$p = $ppt.ActivePresentation

# Locate the slide whose persistent SlideID is 265 (falls back to the
# known position, 11, if that ever changes).
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    if ($p.Slides.Item($i).SlideID -eq 265) {
        $s = $p.Slides.Item($i)
        break
    }
}
if ($s -eq $null) { $s = $p.Slides.Item(11) }

# ---------------------------------------------------------------------
# 1) Shape id=12 ("Text 10"): change run color green -> red
# ---------------------------------------------------------------------
$deltaShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Id -eq 12) {
        $deltaShape = $s.Shapes.Item($i)
        break
    }
}
if ($deltaShape -eq $null) { $deltaShape = $s.Shapes.Item(11) }

$deltaRange = $deltaShape.TextFrame.TextRange
$deltaRange.Font.Color.RGB = 255   # RGB(255,0,0) = FF0000

# ---------------------------------------------------------------------
# 2) Add new shape ("Shape 13") with the conclusion callout
# ---------------------------------------------------------------------
$newShape = $s.Shapes.AddShape(1, 18.0, 371.13, 677.497874015748, 29.57992125984252)
$newShape.Name = "Shape 13"

# White fill, blue outline
$newShape.Fill.ForeColor.RGB = 16777215      # FFFFFF
$newShape.Line.ForeColor.RGB = 16155195      # 3B82F6
$newShape.Line.Weight = 1                    # 12700 EMU = 1 pt
$newShape.Line.DashStyle = 1                 # msoLineSolid

$tf = $newShape.TextFrame
$tr = $tf.TextRange
$tr.Text = "Analizando los resultados el modelo que presenta mejores métricas es Tuning"
$tr.ParagraphFormat.Alignment = 2            # ppAlignCenter

$splitAt = "Analizando los resultados el modelo que presenta mejores métricas es ".Length
$total = $tr.Text.Length

$run2 = $tr.Characters($splitAt + 1, $total - $splitAt)
$run2.LanguageID = "es-MX"
$run2.Font.Size = 24
$run2.Font.Bold = $true

$run1 = $tr.Characters(1, $splitAt)
$run1.LanguageID = "es-MX"
$run1.Font.Size = 20
$run1.Font.Bold = $true

Write-Output "Applied color fix on shape id=$($deltaShape.Id) and added shape id=$($newShape.Id)"
